# Rename the header columns from the generic "_old"/"_new" suffixes to the
# specific format-version suffixes "_FV2310" (old/previous format version)
# and "_FV2404" (new/current format version), add a real Excel Table
# (ListObject) over the used range, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$leftLetters  = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightLetters = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($leftLetters[$i]  + "1").Value = $baseNames[$i] + "_FV2310"
    $ws.Range($rightLetters[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}

# Column K stays "diff" (unchanged).

# Turn the used range into an Excel Table ("Table1") with headers.
$tableRange = $ws.Range("A1:U76")
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true | Out-Null
